# Benchmark update: 2026-01-02 06:43:00 UTC
# Applies the cell-level value changes captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (HESAPTAN EFT - Şube) ---
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = ""
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# --- Row 4 (HESAPTAN EFT - ATM) ---
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = ""
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# --- Row 5 (HESAPTAN EFT - Mobil) ---
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = ""
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# --- Row 6 (DÜZENLİ EFT) ---
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# --- Row 8 (HESAPTAN HAVALE - Şube) ---
$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = ""
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# --- Row 9 (HESAPTAN HAVALE - ATM) ---
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = ""
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# --- Row 10 (HESAPTAN HAVALE - Mobil) ---
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = ""
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# --- Row 11 (DÜZENLİ HAVALE) ---
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

# --- Row 12 (GİDEN SWIFT) ---
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# --- Row 13 (GELEN SWIFT) ---
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# --- Row 14 (GİDEN SWIFT - Mobil) ---
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
$ws.Range("F14").Value = ""
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"
